$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Text edits in paragraph 2 (the "browser / pdf reader" requirements
#    paragraph). These insert extra words/spaces and rename "pdf reader" to
#    "PDF", matching the published diff's resulting paragraph text.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "این سامانه بطور آنلاین و تحت وب مورد استفاده قرار می‌گیرد؛به این منظور به مرورگری که از ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "این سامانه بطور آنلاین و تحت وب مورد استفاده قرار می‌گیرد؛ به این منظور به مرورگری که از ",
    2) | Out-Null

$d.Content.Find.Execute(
    " پشتیبانی کند نیاز دارد.کاربر باید به اینترنت متصل باشد و برای مشاهده ابلاغیه‌ها ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " پشتیبانی کند نیاز دارد. کاربر باید به اینترنت متصل باشد و برای مشاهده ابلاغیه‌ها ",
    2) | Out-Null

$d.Content.Find.Execute(
    "به نرم‌افزار ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "به نرم‌افزاری جهت خواندن فایل های ",
    2) | Out-Null

$d.Content.Find.Execute(
    "pdf reader",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PDF",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Bookmarks
#    - "_Hlk102745234" wraps paragraph 2 (the paragraph just edited above).
#    - "_Hlk102745127" wraps every paragraph from "قیود طراحی" through the
#      final "پشتیبانی مناسب" paragraph.
# ---------------------------------------------------------------------------

$designConstraintsPara = $null
$lastPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("قیود طراحی")) {
        $designConstraintsPara = $i
    }
}
$lastPara = $d.Paragraphs.Count

$browserPara = $d.Paragraphs.Item(2)
$d.Bookmarks.Add("_Hlk102745234", $browserPara.Range) | Out-Null

$startPara = $d.Paragraphs.Item($designConstraintsPara)
$endPara = $d.Paragraphs.Item($lastPara)
$spanRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$d.Bookmarks.Add("_Hlk102745127", $spanRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) Drop the "cs" hint from the rFonts of the last paragraph's paragraph
#    mark run properties (w:ind w:left="360" paragraph / "یا رفع مشکلات...").
# ---------------------------------------------------------------------------

$lastParaObj = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastParaObj.Range.Font.NameBi = "Arial"
